$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (styles/number formats) from the last existing data row
# down into the two new rows before filling in their values.
$ws.Range("A177:F177").Copy()
$ws.Range("A178:F179").PasteSpecial(-4122)

# New row: 2020-09-06
$ws.Range("A178").Value = 44080
$ws.Range("B178").Value = 895
$ws.Range("C178").Value = 412
$ws.Range("D178").Value = 670
$ws.Range("E178").Value = 130
$ws.Range("F178").Value = 30

# New row: 2020-09-07
$ws.Range("A179").Value = 44081
$ws.Range("B179").Value = 281
$ws.Range("C179").Value = 50
$ws.Range("D179").Value = 689
$ws.Range("E179").Value = 123
$ws.Range("F179").Value = 36

# Grow the "Condicion_Pacientes" table so it (and its AutoFilter) cover
# the two newly added rows.
$lo = $ws.ListObjects.Item("Condicion_Pacientes")
$lo.Resize($ws.Range("A1:F179"))

# Match the saved selection from the source workbook.
$ws.Range("F180").Select() | Out-Null
